$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in the data (current data ends at row 93)
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Copy the formatting of the last data row down into the new row so the
# new cells pick up the same style (s="1") as the rest of the table.
$srcRange = $ws.Range("A" + $lastRow + ":E" + $lastRow)
$dstRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$srcRange.Copy($dstRange)

# New stream stat row: stream #92, avg 214, max 262, follow 70,
# games "jc, lol, cs, rdr2, marbels"
$ws.Cells.Item($newRow, 1).Value = 92
$ws.Cells.Item($newRow, 2).Value = 214
$ws.Cells.Item($newRow, 3).Value = 262
$ws.Cells.Item($newRow, 4).Value = 70
$ws.Cells.Item($newRow, 5).Value = "jc, lol, cs, rdr2, marbels"
